$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kareli")

# Copy the formatting from column J (2022) onto the new column K (2023)
# so the new cells pick up the same number formats / fonts / borders
# used throughout the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the 2023 figures
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1066.9000000000001
$ws.Range("K5").Value = 616.29999999999995
$ws.Range("K6").Value = 1319.5
